$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.739.74"
$ws.Range("E2").Value = "  +2.96%  "
$ws.Range("D3").Value = "1.852.15"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.79"
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6373"
$ws.Range("E6").Value = "  +5.23%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07527"
$ws.Range("E8").Value = "  +3.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2976"
$ws.Range("E9").Value = "  +4.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.18"
$ws.Range("E10").Value = "  +6.39%  "
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").Value = "1.840.20"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.057"
$ws.Range("E13").Value = "  +3.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6878"
$ws.Range("E14").Value = "  +5.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "84.53"
$ws.Range("E15").Value = "  +4.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009603"
$ws.Range("E16").Value = "  +7.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.081"
$ws.Range("E17").Value = "  +4.66%  "
$ws.Range("D18").Value = "29.700.50"
$ws.Range("E18").Value = "  +2.87%  "
$ws.Range("D19").Value = "2.094.44"
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "239.54"
$ws.Range("E20").Value = "  +1.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.62"
$ws.Range("E21").Value = "  +2.29%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.356"
$ws.Range("E23").Value = "  +4.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.66"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1421"
$ws.Range("E26").Value = "  +2.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.535"
$ws.Range("E27").Value = "  +2.19%  "
$ws.Range("E28").Value = "  +2.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.501"
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06017"
$ws.Range("E30").Value = "  +8.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.264"
$ws.Range("E31").Value = "  +4.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.146"
$ws.Range("E32").Value = "  +2.48%  "
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.879"
$ws.Range("E34").Value = "  +4.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7331"
$ws.Range("E35").Value = "  +0.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.150"
$ws.Range("E36").Value = "  +2.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.610"
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.861"
$ws.Range("E38").Value = "  +2.11%  "
$ws.Range("D39").Value = "1.226.92"
$ws.Range("E39").Value = "  +3.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01780"
$ws.Range("E40").Value = "  +2.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.358"
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9170"
$ws.Range("E42").Value = "  +3.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "2.008.83"
$ws.Range("E44").Value = "  +2.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.28"
$ws.Range("E45").Value = "  +2.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.45"
$ws.Range("E46").Value = "  +4.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000123"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.336"
$ws.Range("E49").Value = "  +4.12%  "
$ws.Range("E50").Value = "  +3.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1141"
$ws.Range("E51").Value = "  +4.42%  "
